$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: updated prices (D) and hourly volume deltas (E),
# plus a Stellar/Aave row swap (rows 44-45), per the GitHub Actions data pull.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.232.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.578.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.583.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("E12").Value = "  +11.08%  "
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.032.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.244.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.47%  "
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.586.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("E25").Value = "  +7.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.01%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0977"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.953.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
